# "Missing slot for DM added"
#
# - G8: fix the duplicated word typo in the OS-lab slot
#       "آز سیستم عامل (مهندس مهندس شاه منصوری)" -> "آز سیستم عامل (مهندس شاه منصوری)"
#       (the second, purple "ریاضی عمومی 2 (دکتر اعتبار)" run is unchanged)
# - G5: the previously-empty slot now carries the missing Data Mining (DM)
#       course text "مقدمه ای بر داده کاوی (دکتر خیرخواه)", styled the same
#       blue used elsewhere in that column (F5/H5/D5)
# - selection moves onto the newly-populated slot (G4:H4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G8: correct "مهندس مهندس" -> "مهندس" -------------------------------
$g8 = $ws.Range("G8")
$g8NewText = "`nآز سیستم عامل (مهندس شاه منصوری)`nریاضی عمومی 2 (دکتر اعتبار)"
$g8.Value = $g8NewText

# re-apply the purple formatting to the trailing "ریاضی عمومی 2 (دکتر اعتبار)" run
$g8MathLabel = "ریاضی عمومی 2 (دکتر اعتبار)"
$g8MathStart = $g8NewText.Length - $g8MathLabel.Length + 1
$g8MathRun = $g8.Characters($g8MathStart, $g8MathLabel.Length)
$g8MathRun.Font.Color = 10498160   # RGB(0x70,0x30,0xA0) -> BGR long for COM
$g8MathRun.Font.Name = "B Nazanin"
$g8MathRun.Font.Size = 18

# --- G5: fill in the missing "مقدمه ای بر داده کاوی" (Data Mining) slot --
$g5 = $ws.Range("G5")
$g5.Value = "`nمقدمه ای بر داده کاوی (دکتر خیرخواه)"
$g5.Font.Color = 12611584          # RGB(0x00,0x70,0xC0) -> BGR long for COM
$g5.Font.Name = "B Nazanin"
$g5.Font.Size = 18

# --- reflect the newly added slot in the active selection ---------------
$ws.Range("G4:H4").Select()
